$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the custom column width (matching columns B:J) across to column O
$ws.Range("K1:O1").ColumnWidth = 7.83

# Remove the bottom border on the header years row (B3:J3) so only the top border remains
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142  # xlLineStyleNone

# Set the new 2023 values in column K
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 3680
$ws.Range("K5").Value = 1174
$ws.Range("K6").Value = 2506

# Copy column J formatting (font, fill, number format, alignment, borders) into K
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

# Add a right border along the new last column K
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

Write-Host "done"
